# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Most rows just get updated Price (D) / Volume(1h) (E) figures; a few
# rows also swap which coin occupies that rank (name/link/price/volume all
# change together) because the underlying ranking reshuffled.
#
# NumberFormat = "@" is applied before writing any Price value that Excel's
# COM layer would otherwise auto-detect as a pure number (risking silent
# precision loss, e.g. "1.30" collapsing to 1.3) so the text is preserved
# exactly as scraped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.499.85"
$ws.Range("E2").Value = "  -4.55%  "
$ws.Range("D3").Value = "2.491.22"
$ws.Range("E3").Value = "  -5.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.76"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.83"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "2.519.10"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.48"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "2.938.78"
$ws.Range("E14").Value = "  -5.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.53"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").Value = "59.519.44"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "2.499.65"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.65"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.36"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.447"
$ws.Range("E25").Value = "  -10.44%  "
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.82"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.30"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0793"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.83"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.76"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.61"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.97"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("E37").Value = "  -5.23%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.03"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "312.73"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.68"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("E43").Value = "  -7.35%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.07"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0533"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0938"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.64"
$ws.Range("E51").Value = "  -4.80%  "
